$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- View changes: scroll window so column W is the left-most visible
#     column, and move the active selection to AG9 --------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 23   # column W = 23rd column
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AG9").Select()

# --- New column: widen column B (added customWidth column entry) ----
$ws.Columns("B").ColumnWidth = 12.8

# --- Swap the AC:AF values between row 2 and row 3 -------------------
$ws.Range("AC2:AF2").Value = "A"
$ws.Range("AC3:AF3").Value = "I"

Write-Output "edit complete"
